# Applies: "added institut ID and case ID to DQ Reports"
$wb = $excel.ActiveWorkbook

# --- Sheet 1: DQ_Report ---
$ws1 = $wb.Worksheets.Item("DQ_Report")

# Insert a new column before column B ("ICD_primaerkode" etc. shift right)
$ws1.Columns.Item(2).Insert()

# Header for the newly inserted column - copy formatting from a neighboring header cell
$ws1.Range("A1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("B1").Value = "Aufnahmenummer"

# Fill in the new "Aufnahmenummer" values for each data row
$aufnahmenummer = @{
    2  = "F_101645"
    3  = "F_101646"
    4  = "F_101648"
    5  = "F_101649"
    6  = "F_101650"
    7  = "F_101651"
    8  = "F_101651"
    9  = "F_101653"
    10 = "F_101654"
    11 = "F_101655"
    12 = "F_101656"
    13 = "F_101757"
    14 = "F_101658"
    15 = "F_101660"
}

foreach ($row in $aufnahmenummer.Keys) {
    $ws1.Cells.Item($row, 2).Value = $aufnahmenummer[$row]
}

# --- Sheet 2: Statistik ---
$ws2 = $wb.Worksheets.Item("Statistik")
$ws2.Range("A1").Value = "inst_id"
$ws2.Range("A2").Value = "260123430-Dali"
$ws2.Range("B2").Value = 3.13
$ws2.Range("C2").Value = 96.87
